$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell updates reflecting the refreshed cryptocurrency price/volume snapshot.
# Numeric-looking price strings (e.g. "1.000") are forced to Text format first
# so Excel keeps them as literal text instead of coercing them to numbers.

$ws.Range('D2').Value = '27.251.24'
$ws.Range('E2').Value = '  +1.54%  '
$ws.Range('D3').Value = '1.907.20'
$ws.Range('E3').Value = '  +2.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.40'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5247'
$ws.Range('E7').Value = '  +3.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3782'
$ws.Range('E8').Value = '  +3.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07259'
$ws.Range('E9').Value = '  +1.47%  '
$ws.Range('E10').Value = '  +3.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9003'
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08283'
$ws.Range('E12').Value = '  +10.79%  '
$ws.Range('D13').Value = '1.908.85'
$ws.Range('E13').Value = '  +2.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '95.41'
$ws.Range('E14').Value = '  +1.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.281'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008604'
$ws.Range('E17').Value = '  +1.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.49'
$ws.Range('E18').Value = '  +2.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9993'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').Value = '27.282.92'
$ws.Range('E20').Value = '  +1.55%  '
$ws.Range('E21').Value = '  +1.53%  '
$ws.Range('D22').Value = '2.152.52'
$ws.Range('E22').Value = '  +1.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.66'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.458'
$ws.Range('E24').Value = '  +1.78%  '
$ws.Range('E25').Value = '  +10.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '145.96'
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.743'
$ws.Range('E27').Value = '  -1.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.16'
$ws.Range('E28').Value = '  +1.76%  '
$ws.Range('E29').Value = '  +1.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.974'
$ws.Range('E30').Value = '  +5.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.815'
$ws.Range('E31').Value = '  +3.14%  '
$ws.Range('E32').Value = '  +0.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8039'
$ws.Range('E33').Value = '  +7.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05085'
$ws.Range('E34').Value = '  +1.20%  '
$ws.Range('E35').Value = '  +7.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.936'
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.355'
$ws.Range('E37').Value = '  +4.60%  '
$ws.Range('E38').Value = '  +2.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5742'
$ws.Range('E39').Value = '  +3.51%  '
$ws.Range('E40').Value = '  +0.18%  '
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.070'
$ws.Range('E42').Value = '  +5.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.639'
$ws.Range('E43').Value = '  +1.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '118.44'
$ws.Range('E44').Value = '  +2.22%  '
$ws.Range('E45').Value = '  +2.40%  '
$ws.Range('E46').Value = '  +2.29%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.000'
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.15'
$ws.Range('E48').Value = '  +2.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.614'
$ws.Range('E49').Value = '  +4.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.61'
$ws.Range('E50').Value = '  +1.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.70'
$ws.Range('E51').Value = '  +1.53%  '
